$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "65.647.71"
$ws.Range("E2").Value = "  +3.26%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.500.66"
$ws.Range("E3").Value = "  +2.58%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.03%  "

# Row 5 - BNB
Set-TextValue "D5" "581.12"
$ws.Range("E5").Value = "  +2.04%  "

# Row 6 - Solana
Set-TextValue "D6" "161.43"
$ws.Range("E6").Value = "  +3.92%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.06%  "

# Row 8 - XRP
Set-TextValue "D8" "0.608"
$ws.Range("E8").Value = "  +11.72%  "

# Row 9 - LidoStakedEther
Set-TextValue "D9" "3.503.23"
$ws.Range("E9").Value = "  +2.50%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  -1.36%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +2.67%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +2.25%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "4.105.23"
$ws.Range("E13").Value = "  +2.21%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -0.03%  "

# Row 15 - ShibaInu
$ws.Range("E15").Value = "  +2.75%  "

# Row 16 - Avalanche
Set-TextValue "D16" "28.65"
$ws.Range("E16").Value = "  +5.92%  "

# Row 17 - now WrappedEther (was WrappedBTC)
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D17" "3.579.56"
$ws.Range("E17").Value = "  +4.33%  "

# Row 18 - now WrappedBTC (was WrappedEther)
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D18" "65.666.38"
$ws.Range("E18").Value = "  +3.09%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +2.44%  "

# Row 20 - Chainlink
Set-TextValue "D20" "14.29"
$ws.Range("E20").Value = "  +0.76%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "390.27"
$ws.Range("E21").Value = "  +1.16%  "

# Row 22 - Uniswap
Set-TextValue "D22" "8.28"
$ws.Range("E22").Value = "  +1.01%  "

# Row 23 - Polygon
Set-TextValue "D23" "0.551"
$ws.Range("E23").Value = "  +2.71%  "

# Row 24 - Litecoin
Set-TextValue "D24" "73.46"
$ws.Range("E24").Value = "  +1.55%  "

# Row 25 - Dai
$ws.Range("E25").Value = "  +0.71%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  +6.03%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue "D27" "9.99"
$ws.Range("E27").Value = "  +5.37%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  +1.69%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  -0.10%  "

# Row 30 - NEARProtocol
$ws.Range("E30").Value = "  +6.80%  "

# Row 31 - Fetch.AI
$ws.Range("E31").Value = "  +7.98%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  +2.82%  "

# Row 33 - EthereumClassic
$ws.Range("E33").Value = "  +1.68%  "

# Row 34 - RenderToken
Set-TextValue "D34" "6.51"
$ws.Range("E34").Value = "  +1.24%  "

# Row 36 - Aptos
Set-TextValue "D36" "7.17"
$ws.Range("E36").Value = "  +3.62%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  +7.23%  "

# Row 38 - Monero
Set-TextValue "D38" "162.67"
$ws.Range("E38").Value = "  +1.50%  "

# Row 39 - Stacks
Set-TextValue "D39" "1.95"
$ws.Range("E39").Value = "  +7.19%  "

# Row 40 - Maker
Set-TextValue "D40" "3.092.50"
$ws.Range("E40").Value = "  +7.12%  "

# Row 41 - Hedera
Set-TextValue "D41" "0.0773"
$ws.Range("E41").Value = "  -0.12%  "

# Row 42 - EnergySwap
$ws.Range("E42").Value = "  +1.97%  "

# Row 43 - VeChain
$ws.Range("E43").Value = "  +2.03%  "

# Row 44 - Filecoin
$ws.Range("E44").Value = "  +4.07%  "

# Row 45 - OKB
Set-TextValue "D45" "43.22"
$ws.Range("E45").Value = "  +4.57%  "

# Row 46 - Mantle
Set-TextValue "D46" "0.784"
$ws.Range("E46").Value = "  +2.64%  "

# Row 47 - InjectiveProtocol
Set-TextValue "D47" "26.01"
$ws.Range("E47").Value = "  +11.47%  "

# Row 48 - ONDO
$ws.Range("E48").Value = "  +5.07%  "

# Row 49 - dogwifhat
Set-TextValue "D49" "2.25"
$ws.Range("E49").Value = "  +4.57%  "

# Row 50 - Bittensor
Set-TextValue "D50" "315.05"
$ws.Range("E50").Value = "  +8.89%  "

# Row 51 - Cosmos
Set-TextValue "D51" "6.73"
$ws.Range("E51").Value = "  +4.12%  "
